# Blue Mountain Community College Organizations workbook edit
# - Swap "Organization Name" / "Category" columns (A<->B), renaming headers
# - Rename several other headers
# - Add new "Tiktok Link" column (M)
# - Adjust column widths to match new layout

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 15

# 1) Swap the data in columns A and B for every data row (2..15).
#    Row 1 (headers) is handled separately below since the text itself
#    changes, not just the position.
for ($r = 2; $r -le $lastRow; $r++) {
    $aVal = $ws.Cells.Item($r, 1).Value2
    $bVal = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($r, 1).Value = $bVal
    $ws.Cells.Item($r, 2).Value = $aVal
}

# 2) Update header row text (row 1)
$ws.Range("A1").Value = "Category"
$ws.Range("B1").Value = "Organization Name"
$ws.Range("C1").Value = "Organization Link"
$ws.Range("D1").Value = "Logo Link"
$ws.Range("E1").Value = "Description"
$ws.Range("F1").Value = "Email"
$ws.Range("G1").Value = "Phone Number"
$ws.Range("H1").Value = "Linkedin Link"
$ws.Range("I1").Value = "Instagram Link"
$ws.Range("J1").Value = "Facebook Link"
$ws.Range("K1").Value = "Twitter Link"
$ws.Range("L1").Value = "Youtube Link"

# M1 is a brand new header cell; copy the header style (bold font, border,
# centered alignment, etc.) from an existing header cell before setting its
# text so it matches the rest of row 1 (s="1" in the saved XML).
$ws.Range("L1").Copy()
$ws.Range("M1").PasteSpecial(-4122)
$ws.Range("M1").Value = "Tiktok Link"

# 3) Column widths.
#    Excel's Range/Columns.ColumnWidth is measured in "characters" and is
#    offset from the raw OOXML <col width> by the standard padding
#    constant (5/6 = 0.8333...). Subtract that here so the saved XML ends
#    up with the exact widths we want.
$widthOffset = 0.8333333333333334

$ws.Columns.Item(1).ColumnWidth  = 20 - $widthOffset
$ws.Columns.Item(2).ColumnWidth  = 35 - $widthOffset
$ws.Columns.Item(7).ColumnWidth  = 14 - $widthOffset
$ws.Columns.Item(8).ColumnWidth  = 15 - $widthOffset
$ws.Columns.Item(9).ColumnWidth  = 16 - $widthOffset
$ws.Columns.Item(10).ColumnWidth = 15 - $widthOffset
$ws.Columns.Item(11).ColumnWidth = 14 - $widthOffset
$ws.Columns.Item(12).ColumnWidth = 14 - $widthOffset
$ws.Columns.Item(13).ColumnWidth = 13 - $widthOffset
